# Update the "more todo:" list on Sheet1:
#  - reorder existing items ("test all footprints", "silkscreen") up to
#    rows 21/22, dropping the stray font-only style on row 22
#  - replace the two finished/old TODO items (now rows 23/24) with the
#    new hoekstep-related TODO text
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A21").Value = "test all footprints"

$ws1.Range("A22").Value = "silkscreen"
$ws1.Range("A22").ClearFormats() | Out-Null

$ws1.Range("A23").Value = "add config jumpers for hoekstep drivers"

$ws1.Range("A24").Value = "look into routing DAC into ROSC pins on hoekstep"

# Sheet2's selection (C15) is untouched by this edit; re-assert it first.
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate() | Out-Null
$ws2.Range("C15").Select() | Out-Null

# Sheet1 stays the active/visible tab; move the active selection on it
# down to the last edited cell.
$ws1.Activate() | Out-Null
$ws1.Range("A24").Select() | Out-Null
